# Add a new row of "login" data (username, fullname, password) to the
# database sheet, matching the pattern of the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "hello21"
$ws.Range("B6").Value = "hello"

# C6 must stay a text value ("1234"), not be auto-coerced to a number the
# way Excel normally would. Force text via NumberFormat, assign the value,
# then clear the formatting override so the cell ends up with the same
# (default) style as its neighbours while keeping its text type.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1234"
$ws.Range("C6").ClearFormats()
